# Update gh-pages output: new event inserted into the "展览" (exhibitions)
# sheet (and the merged "全部类型" sheet), plus "想去人数" (interest count)
# refreshes across both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (index 1)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Refresh "想去人数" (F column) counts for rows that are not shifting.
$ws1.Range("F4").Value = 78
$ws1.Range("F5").Value = 4
$ws1.Range("F6").Value = 542
$ws1.Range("F7").Value = 7563
$ws1.Range("F8").Value = 479
$ws1.Range("F10").Value = 1070
$ws1.Range("F11").Value = 594

# Make room for the new event by shifting rows 12-16 down to 13-17,
# copying bottom-up so nothing is overwritten before it is copied.
for ($r = 16; $r -ge 12; $r--) {
    $ws1.Rows.Item($r).Copy($ws1.Rows.Item($r + 1))
}

# Write the newly-added event into row 12.
$ws1.Range("B12").Value = "2024-06-30"
$ws1.Range("C12").Value = "合肥·第1.5届星芒动漫嘉年华"
$ws1.Range("D12").Value = "山西路与太原路交叉口 挥动体育"
$ws1.Range("E12").Value = "2024.06.30 09:30-06.30 17:30"
$ws1.Range("F12").Value = 3
$ws1.Range("G12").Value = 60
$ws1.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=85213"
$ws1.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202405/v40vLtJl1715073148563.jpeg"

# The events that were shifted down also got their interest counts bumped.
$ws1.Range("F13").Value = 19
$ws1.Range("F14").Value = 169
$ws1.Range("F15").Value = 2
$ws1.Range("F16").Value = 198
$ws1.Range("F17").Value = 723

# ---------------------------------------------------------------------
# Sheet "全部类型" (index 4) - mirrors "展览" plus the "演出" rows
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 78
$ws4.Range("F5").Value = 4
$ws4.Range("F7").Value = 542
$ws4.Range("F8").Value = 7563
$ws4.Range("F9").Value = 479
$ws4.Range("F11").Value = 1070
$ws4.Range("F12").Value = 594

# Shift rows 13-18 down to 14-19 for the new row, bottom-up.
for ($r = 18; $r -ge 13; $r--) {
    $ws4.Rows.Item($r).Copy($ws4.Rows.Item($r + 1))
}

$ws4.Range("B13").Value = "2024-06-30"
$ws4.Range("C13").Value = "合肥·第1.5届星芒动漫嘉年华"
$ws4.Range("D13").Value = "山西路与太原路交叉口 挥动体育"
$ws4.Range("E13").Value = "2024.06.30 09:30-06.30 17:30"
$ws4.Range("F13").Value = 3
$ws4.Range("G13").Value = 60
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=85213"
$ws4.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202405/v40vLtJl1715073148563.jpeg"

$ws4.Range("F14").Value = 19
$ws4.Range("F15").Value = 169
$ws4.Range("F16").Value = 2
$ws4.Range("F17").Value = 198
$ws4.Range("F18").Value = 723
